$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the new columns I (I0) and J (IF). Copy the formatting
# (bold font, border, centered alignment) from the existing header cell
# H1 so the new headers reuse the same cell style as the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill data rows 2-36: column I is always 1, column J mirrors column H.
for ($r = 2; $r -le 36; $r++) {
    $h = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $h
}
